$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.214864207225276
$ws.Cells.Item(2, 3).Value = 0.2001994454756755
$ws.Cells.Item(2, 4).Value = 0.5023494010699778
$ws.Cells.Item(2, 5).Value = 0.1702098314471225
$ws.Cells.Item(2, 7).Value = 0.8783912334407518
$ws.Cells.Item(2, 8).Value = 0.9460361915521531
$ws.Cells.Item(2, 9).Value = 0.9225269522909834
$ws.Cells.Item(2, 10).Value = 0.07982560351849344
$ws.Cells.Item(2, 12).Value = 0.4054404887654783
$ws.Cells.Item(2, 13).Value = 0.3290411884122832
$ws.Cells.Item(2, 15).Value = 3.680294265256407

# Row 3
$ws.Cells.Item(3, 2).Value = 1.110751157864399
$ws.Cells.Item(3, 3).Value = 0.1861287943088143
$ws.Cells.Item(3, 4).Value = 0.5035493356212584
$ws.Cells.Item(3, 5).Value = 0.1716971708228554
$ws.Cells.Item(3, 7).Value = 0.8858018408170025
$ws.Cells.Item(3, 8).Value = 0.9542607946025257
$ws.Cells.Item(3, 9).Value = 0.9351310122851899
$ws.Cells.Item(3, 10).Value = 0.07923726039927104
$ws.Cells.Item(3, 12).Value = 0.4012443182651708
$ws.Cells.Item(3, 13).Value = 0.3113875287910304
$ws.Cells.Item(3, 15).Value = 3.712910337738762

# Row 4
$ws.Cells.Item(4, 2).Value = 1.046814342849757
$ws.Cells.Item(4, 3).Value = 0.1774485326391755
$ws.Cells.Item(4, 4).Value = 0.504503896581511
$ws.Cells.Item(4, 5).Value = 0.1726667698276714
$ws.Cells.Item(4, 7).Value = 0.8909588459023752
$ws.Cells.Item(4, 8).Value = 0.9597548147556765
$ws.Cells.Item(4, 9).Value = 0.9433980698569684
$ws.Cells.Item(4, 10).Value = 0.07887971963609175
$ws.Cells.Item(4, 12).Value = 0.3988052554240227
$ws.Cells.Item(4, 13).Value = 0.3005991828378143
$ws.Cells.Item(4, 15).Value = 3.735141989268811

# Row 5
$ws.Cells.Item(5, 2).Value = 1.020758700332749
$ws.Cells.Item(5, 3).Value = 0.1739012395003101
$ws.Cells.Item(5, 4).Value = 0.5049477411495502
$ws.Cells.Item(5, 5).Value = 0.1730760808267489
$ws.Cells.Item(5, 7).Value = 0.8932128530664301
$ws.Cells.Item(5, 8).Value = 0.9621053889245132
$ws.Cells.Item(5, 9).Value = 0.9468997088524347
$ws.Cells.Item(5, 10).Value = 0.07873496897397914
$ws.Cells.Item(5, 12).Value = 0.397846014210117
$ws.Cells.Item(5, 13).Value = 0.2962160285197868
$ws.Cells.Item(5, 15).Value = 3.744755947444588

# Row 6
$ws.Cells.Item(6, 2).Value = 1.01643218399056
$ws.Cells.Item(6, 3).Value = 0.1733116166213335
$ws.Cells.Item(6, 4).Value = 0.5050247567479431
$ws.Cells.Item(6, 5).Value = 0.173144904259674
$ws.Cells.Item(6, 7).Value = 0.8935963356291268
$ws.Cells.Item(6, 8).Value = 0.9625024487999951
$ws.Cells.Item(6, 9).Value = 0.9474891693228553
$ws.Cells.Item(6, 10).Value = 0.07871099117103952
$ws.Cells.Item(6, 12).Value = 0.3976888328266739
$ws.Cells.Item(6, 13).Value = 0.2954890154736489
$ws.Cells.Item(6, 15).Value = 3.746385816678355

# Row 7
$ws.Cells.Item(7, 2).Value = 1.046462948356293
$ws.Cells.Item(7, 3).Value = 0.1774007328256744
$ws.Cells.Item(7, 4).Value = 0.5045096602103172
$ws.Cells.Item(7, 5).Value = 0.1726722324543957
$ws.Cells.Item(7, 7).Value = 0.8909886269757408
$ws.Cells.Item(7, 8).Value = 0.9597860630219657
$ws.Cells.Item(7, 9).Value = 0.9434447567726139
$ws.Cells.Item(7, 10).Value = 0.07887776360627186
$ws.Cells.Item(7, 12).Value = 0.3987921780819761
$ws.Cells.Item(7, 13).Value = 0.3005400162488527
$ws.Cells.Item(7, 15).Value = 3.7352694020108

# Row 8
$ws.Cells.Item(8, 2).Value = 1.17896940495848
$ws.Cells.Item(8, 3).Value = 0.1953565032417544
$ws.Cells.Item(8, 4).Value = 0.5027179860282303
$ws.Cells.Item(8, 5).Value = 0.1707109779148164
$ws.Cells.Item(8, 7).Value = 0.8808203872069811
$ws.Cells.Item(8, 8).Value = 0.9487798950606177
$ws.Cells.Item(8, 9).Value = 0.9267631916843193
$ws.Cells.Item(8, 10).Value = 0.07962198491317451
$ws.Cells.Item(8, 12).Value = 0.4039652164130558
$ws.Cells.Item(8, 13).Value = 0.3229438051350613
$ws.Cells.Item(8, 15).Value = 3.691082496955474

# Row 9
$ws.Cells.Item(9, 2).Value = 1.438657109959081
$ws.Cells.Item(9, 3).Value = 0.230234508342221
$ws.Cells.Item(9, 4).Value = 0.5009291916411343
$ws.Cells.Item(9, 5).Value = 0.1673113817869192
$ws.Cells.Item(9, 7).Value = 0.8657012766164343
$ws.Cells.Item(9, 8).Value = 0.9307183371165308
$ws.Cells.Item(9, 9).Value = 0.8982430953031049
$ws.Cells.Item(9, 10).Value = 0.08111005888186895
$ws.Cells.Item(9, 12).Value = 0.4151942069556895
$ws.Cells.Item(9, 13).Value = 0.3672701315979126
$ws.Cells.Item(9, 15).Value = 3.621937643179251

# Row 10
$ws.Cells.Item(10, 2).Value = 1.629279268799053
$ws.Cells.Item(10, 3).Value = 0.2556462687513203
$ws.Cells.Item(10, 4).Value = 0.5006622773241389
$ws.Cells.Item(10, 5).Value = 0.1650846650083402
$ws.Cells.Item(10, 7).Value = 0.8575405022835128
$ws.Cells.Item(10, 8).Value = 0.9195926168896733
$ws.Cells.Item(10, 9).Value = 0.879847576787121
$ws.Cells.Item(10, 10).Value = 0.08221995609919119
$ws.Cells.Item(10, 12).Value = 0.4240991817468966
$ws.Cells.Item(10, 13).Value = 0.400062413039187
$ws.Cells.Item(10, 15).Value = 3.581821309816405

# Row 11
$ws.Cells.Item(11, 2).Value = 1.715946749252794
$ws.Cells.Item(11, 3).Value = 0.2671585253259536
$ws.Cells.Item(11, 4).Value = 0.5007674347790072
$ws.Cells.Item(11, 5).Value = 0.1641302592984175
$ws.Cells.Item(11, 7).Value = 0.8544697524940119
$ws.Cells.Item(11, 8).Value = 0.9149962620795122
$ws.Cells.Item(11, 9).Value = 0.8720349732240464
$ws.Cells.Item(11, 10).Value = 0.0827283106323442
$ws.Cells.Item(11, 12).Value = 0.4282912885865073
$ws.Cells.Item(11, 13).Value = 0.4150268424488601
$ws.Cells.Item(11, 15).Value = 3.565894202070979

# Row 12
$ws.Cells.Item(12, 2).Value = 1.748757019383845
$ws.Cells.Item(12, 3).Value = 0.2715108466990728
$ws.Cells.Item(12, 4).Value = 0.500839758929331
$ws.Cells.Item(12, 5).Value = 0.1637772504287458
$ws.Cells.Item(12, 7).Value = 0.8533993451251973
$ws.Cells.Item(12, 8).Value = 0.9133225389350486
$ws.Cells.Item(12, 9).Value = 0.8691564997915151
$ws.Cells.Item(12, 10).Value = 0.08292129152484407
$ws.Cells.Item(12, 12).Value = 0.4298989052470716
$ws.Cells.Item(12, 13).Value = 0.4206999550307202
$ws.Cells.Item(12, 15).Value = 3.560197145911985

# Row 13
$ws.Cells.Item(13, 2).Value = 1.741691164705173
$ws.Cells.Item(13, 3).Value = 0.2705738172769543
$ws.Cells.Item(13, 4).Value = 0.5008227382838442
$ws.Cells.Item(13, 5).Value = 0.1638529036626393
$ws.Cells.Item(13, 7).Value = 0.8536257637148594
$ws.Cells.Item(13, 8).Value = 0.9136800337844591
$ws.Cells.Item(13, 9).Value = 0.8697728714778918
$ws.Cells.Item(13, 10).Value = 0.08287970867614192
$ws.Cells.Item(13, 12).Value = 0.4295517820770698
$ws.Cells.Item(13, 13).Value = 0.4194778694253856
$ws.Cells.Item(13, 15).Value = 3.561409240919403

# Row 14
$ws.Cells.Item(14, 2).Value = 1.718646258137824
$ws.Cells.Item(14, 3).Value = 0.2675167376315812
$ws.Cells.Item(14, 4).Value = 0.5007727339642969
$ws.Cells.Item(14, 5).Value = 0.1641010486999912
$ws.Cells.Item(14, 7).Value = 0.8543798363472348
$ws.Cells.Item(14, 8).Value = 0.9148572248763855
$ws.Cells.Item(14, 9).Value = 0.8717965558265632
$ws.Cells.Item(14, 10).Value = 0.08274417781390042
$ws.Cells.Item(14, 12).Value = 0.4284231450396589
$ws.Cells.Item(14, 13).Value = 0.4154934464577806
$ws.Cells.Item(14, 15).Value = 3.565418802431111

# Row 15
$ws.Cells.Item(15, 2).Value = 1.704529379393364
$ws.Cells.Item(15, 3).Value = 0.2656432530037307
$ws.Cells.Item(15, 4).Value = 0.5007463354385351
$ws.Cells.Item(15, 5).Value = 0.164254138718098
$ws.Cells.Item(15, 7).Value = 0.8548537678187103
$ws.Cells.Item(15, 8).Value = 0.9155869893294977
$ws.Cells.Item(15, 9).Value = 0.8730465407644701
$ws.Cells.Item(15, 10).Value = 0.08266122294411815
$ws.Cells.Item(15, 12).Value = 0.4277344429480507
$ws.Cells.Item(15, 13).Value = 0.4130536962352167
$ws.Cells.Item(15, 15).Value = 3.567918305132963

# Row 16
$ws.Cells.Item(16, 2).Value = 1.623614218632895
$ws.Cells.Item(16, 3).Value = 0.2548929336419121
$ws.Cells.Item(16, 4).Value = 0.5006599576704787
$ws.Cells.Item(16, 5).Value = 0.1651482144229144
$ws.Cells.Item(16, 7).Value = 0.8577541073572093
$ws.Cells.Item(16, 8).Value = 0.9199023505829018
$ws.Cells.Item(16, 9).Value = 0.8803693371315546
$ws.Cells.Item(16, 10).Value = 0.08218680186442739
$ws.Cells.Item(16, 12).Value = 0.4238280462956823
$ws.Cells.Item(16, 13).Value = 0.3990853685646272
$ws.Cells.Item(16, 15).Value = 3.582908933565307

# Row 17
$ws.Cells.Item(17, 2).Value = 1.573961800495908
$ws.Cells.Item(17, 3).Value = 0.2482855602804079
$ws.Cells.Item(17, 4).Value = 0.500664944050925
$ws.Cells.Item(17, 5).Value = 0.165711684165788
$ws.Cells.Item(17, 7).Value = 0.8596978176832692
$ws.Cells.Item(17, 8).Value = 0.9226687150245283
$ws.Cells.Item(17, 9).Value = 0.8850040074696679
$ws.Cells.Item(17, 10).Value = 0.08189663167321726
$ws.Cells.Item(17, 12).Value = 0.4214676563169348
$ws.Cells.Item(17, 13).Value = 0.3905280482982221
$ws.Cells.Item(17, 15).Value = 3.592700129975498

# Row 18
$ws.Cells.Item(18, 2).Value = 1.545398690640582
$ws.Cells.Item(18, 3).Value = 0.2444807015894526
$ws.Cells.Item(18, 4).Value = 0.5006891433771017
$ws.Cells.Item(18, 5).Value = 0.1660412878484889
$ws.Cells.Item(18, 7).Value = 0.8608761723566118
$ws.Cells.Item(18, 8).Value = 0.9243036041933266
$ws.Cells.Item(18, 9).Value = 0.8877220383070394
$ws.Cells.Item(18, 10).Value = 0.08173005994304106
$ws.Cells.Item(18, 12).Value = 0.4201233205191528
$ws.Cells.Item(18, 13).Value = 0.3856105544634119
$ws.Cells.Item(18, 15).Value = 3.598550308839094

# Row 19
$ws.Cells.Item(19, 2).Value = 1.535727025912195
$ws.Cells.Item(19, 3).Value = 0.2431916817695026
$ws.Cells.Item(19, 4).Value = 0.5007010029102901
$ws.Cells.Item(19, 5).Value = 0.166153832902709
$ws.Cells.Item(19, 7).Value = 0.8612855092918466
$ws.Cells.Item(19, 8).Value = 0.924864664199518
$ws.Cells.Item(19, 9).Value = 0.8886512950205692
$ws.Cells.Item(19, 10).Value = 0.08167371827403613
$ws.Cells.Item(19, 12).Value = 0.4196704399348619
$ws.Cells.Item(19, 13).Value = 0.3839463508777854
$ws.Cells.Item(19, 15).Value = 3.600568604713942

# Row 20
$ws.Cells.Item(20, 2).Value = 1.579247850093736
$ws.Cells.Item(20, 3).Value = 0.2489893912391778
$ws.Cells.Item(20, 4).Value = 0.5006622061016373
$ws.Cells.Item(20, 5).Value = 0.1656511316051104
$ws.Cells.Item(20, 7).Value = 0.8594846557749207
$ws.Cells.Item(20, 8).Value = 0.9223697030390099
$ws.Cells.Item(20, 9).Value = 0.8845052262915303
$ws.Cells.Item(20, 10).Value = 0.08192748713322118
$ws.Cells.Item(20, 12).Value = 0.4217175486176643
$ws.Cells.Item(20, 13).Value = 0.3914385316491433
$ws.Cells.Item(20, 15).Value = 3.591635220038853

# Row 21
$ws.Cells.Item(21, 2).Value = 1.725415361111857
$ws.Cells.Item(21, 3).Value = 0.2684148713656498
$ws.Cells.Item(21, 4).Value = 0.5007865399395968
$ws.Cells.Item(21, 5).Value = 0.1640279345199573
$ws.Cells.Item(21, 7).Value = 0.8541558374947442
$ws.Cells.Item(21, 8).Value = 0.9145096420886034
$ws.Cells.Item(21, 9).Value = 0.8711999790255689
$ws.Cells.Item(21, 10).Value = 0.08278397368936652
$ws.Cells.Item(21, 12).Value = 0.4287541072080785
$ws.Cells.Item(21, 13).Value = 0.4166635970057655
$ws.Cells.Item(21, 15).Value = 3.564232024788993

# Row 22
$ws.Cells.Item(22, 2).Value = 1.820891931652795
$ws.Cells.Item(22, 3).Value = 0.2810689217827189
$ws.Cells.Item(22, 4).Value = 0.5010572170516383
$ws.Cells.Item(22, 5).Value = 0.1630160590964183
$ws.Cells.Item(22, 7).Value = 0.8512119053213922
$ws.Cells.Item(22, 8).Value = 0.9097620874813828
$ws.Cells.Item(22, 9).Value = 0.8629705170941655
$ws.Cells.Item(22, 10).Value = 0.08334651497239776
$ws.Cells.Item(22, 12).Value = 0.4334703267659847
$ws.Cells.Item(22, 13).Value = 0.4331868208016303
$ws.Cells.Item(22, 15).Value = 3.548270523450981

# Row 23
$ws.Cells.Item(23, 2).Value = 1.769939734220998
$ws.Cells.Item(23, 3).Value = 0.2743191121547852
$ws.Cells.Item(23, 4).Value = 0.5008954457349546
$ws.Cells.Item(23, 5).Value = 0.1635516391818239
$ws.Cells.Item(23, 7).Value = 0.8527337904509125
$ws.Cells.Item(23, 8).Value = 0.9122603185518727
$ws.Cells.Item(23, 9).Value = 0.8673200395835323
$ws.Cells.Item(23, 10).Value = 0.08304602817398177
$ws.Cells.Item(23, 12).Value = 0.4309424956458088
$ws.Cells.Item(23, 13).Value = 0.4243647808002677
$ws.Cells.Item(23, 15).Value = 3.556611133004424

# Row 24
$ws.Cells.Item(24, 2).Value = 1.576858080724037
$ws.Cells.Item(24, 3).Value = 0.2486712085309648
$ws.Cells.Item(24, 4).Value = 0.5006633774747797
$ws.Cells.Item(24, 5).Value = 0.1656784897944172
$ws.Cells.Item(24, 7).Value = 0.8595808366175532
$ws.Cells.Item(24, 8).Value = 0.9225047478352337
$ws.Cells.Item(24, 9).Value = 0.8847305586280214
$ws.Cells.Item(24, 10).Value = 0.08191353659618983
$ws.Cells.Item(24, 12).Value = 0.4216045327890896
$ws.Cells.Item(24, 13).Value = 0.391026895173141
$ws.Cells.Item(24, 15).Value = 3.59211597715958

# Row 25
$ws.Cells.Item(25, 2).Value = 1.368429290153188
$ws.Cells.Item(25, 3).Value = 0.2208357958250247
$ws.Cells.Item(25, 4).Value = 0.5012288931249316
$ws.Cells.Item(25, 5).Value = 0.1681833930092385
$ws.Cells.Item(25, 7).Value = 0.8692744305605373
$ws.Cells.Item(25, 8).Value = 0.9352277453513835
$ws.Cells.Item(25, 9).Value = 0.905509463819989
$ws.Cells.Item(25, 10).Value = 0.08070451705166093
$ws.Cells.Item(25, 12).Value = 0.4120409897699204
$ws.Cells.Item(25, 13).Value = 0.3552380721899553
$ws.Cells.Item(25, 15).Value = 3.638767765219313

Write-Output "pl_mw case 4_15 (380 kV) values updated"